$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert two new rows above row 2 (existing data / formulas shift down and
# formula ranges auto-adjust).
$ws.Range("A2:A3").EntireRow.Insert()

# Match the date-formatted style used by the rest of column A (copy
# formats only, so no new style/numFmt entries are introduced).
$ws.Range("A4").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data for 18 and 19 May (row 2 = most recent date, 45431).
$ws.Range("A2").Value = 45431
$ws.Range("B2").Value = 240
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 211

$ws.Range("A3").Value = 45430
$ws.Range("B3").Value = 213
$ws.Range("C3").Value = 23
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 187

# Remove the stray column I annotations (now unused / out of range).
$ws.Range("I1:I21").ClearContents()

# Extend the totals formulas to cover the newly-added rows.
$ws.Range("B21").Formula = "=SUM(B2:B20)"
$ws.Range("C21:E21").Formula = "=SUM(C2:C20)"

$ws.Range("B21:E21").Select()
$excel.ActiveWindow.ScrollRow = 2
